# Update "想去人数" (number of people interested) counts in the 漫展信息 workbook
# to match the freshly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1232   # 广州·南国书香节璃樱动漫嘉年华展区【免费入场】
$wsExpo.Range("F4").Value = 12804  # 广州·第九届初物语动漫展
$wsExpo.Range("F5").Value = 727    # 广州·COC星火次元云漫创作交流展
$wsExpo.Range("F16").Value = 345   # 广州·第五人格only同人展
$wsExpo.Range("F18").Value = 296   # 广州·少女番only4.0
$wsExpo.Range("F23").Value = 246   # 广州·wio流金序曲乙女同人展
$wsExpo.Range("F24").Value = 1287  # 广州·第一届Redamancy动漫游戏嘉年华

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F6").Value = 156    # 广州·春日计划2024——特别二次元不插电音乐会

# Sheet "本地生活" (Local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 4025  # 广州·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1232    # 广州·南国书香节璃樱动漫嘉年华展区【免费入场】
$wsAll.Range("F7").Value = 12804   # 广州·第九届初物语动漫展
$wsAll.Range("F9").Value = 727     # 广州·COC星火次元云漫创作交流展
$wsAll.Range("F10").Value = 4025   # 广州·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅
$wsAll.Range("F21").Value = 156    # 广州·春日计划2024——特别二次元不插电音乐会
$wsAll.Range("F22").Value = 156    # 广州·春日计划2024——特别二次元不插电音乐会
$wsAll.Range("F29").Value = 345    # 广州·第五人格only同人展
$wsAll.Range("F32").Value = 296    # 广州·少女番only4.0
$wsAll.Range("F40").Value = 246    # 广州·wio流金序曲乙女同人展
$wsAll.Range("F41").Value = 1287   # 广州·第一届Redamancy动漫游戏嘉年华
